# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holdings detail) right before
#    the "总计" (totals) summary sheet.
# 2. Populate it with the Q1-2022 fund holding rows, mirroring the layout
#    used by the other quarterly sheets (e.g. "2021-Q4").
# 3. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locate the existing sheets we need as templates / anchors.
# ---------------------------------------------------------------------
$totalSheetAnchor = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet directly before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetAnchor)
$newSheet.Name = "2022-Q1"

# Borrow the header-row (B1:H1) and "index column" (A) formatting from an
# existing quarterly sheet so styles match (bold font + border + center
# alignment on headers, centered bold on column A).
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)
$wb.Application.CutCopyMode = 0

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比,
# 持有市值(亿元), 仓位排名 -- columns B..H, with a 0-based row index in A.
$rows = @(
    @("011913", "华夏永泓一年持有混合A", "24.65", "37.51", "1.05", "0.2588", 6),
    @("011914", "华夏永泓一年持有混合C", "12.39", "37.51", "1.05", "0.1301", 6),
    @("513690", "博时恒生港股通高股息率ETF", "4.60", "99.64", "2.55", "0.1173", 6),
    @("011355", "华泰柏瑞港股通时代机遇混合型证券投资基金A", "1.13", "90.93", "6.47", "0.0731", 4),
    @("003413", "华泰柏瑞新经济沪港深灵活配置混合", "0.54", "92.57", "5.39", "0.0291", 10),
    @("011356", "华泰柏瑞港股通时代机遇混合型证券投资基金C", "0.40", "90.93", "6.47", "0.0259", 4),
    @("005702", "恒生前海港股通高股息低波动指数", "0.29", "94.14", "2.60", "0.0075", 3)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = "'" + $data[0]
    $newSheet.Cells.Item($r, 3).Value = $data[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $data[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $data[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $data[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $data[5]
    $newSheet.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Prepend the "2022-Q1" summary row to "总计", shifting rows 2.. down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$wb.Application.CutCopyMode = 0

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.64
